# Updated symbol list on Tue Dec 13 17:57:49 UTC 2022 with GitHub Actions
#
# Applies the refreshed crypto price / volume snapshot to Sheet1.
# Every cell below is stored as text in the workbook (e.g. "268.43"),
# so for the numeric-looking values we explicitly switch the cell to
# Text format before writing, otherwise Excel would silently convert
# the string into a floating point number (and drop meaningful
# trailing zeros, e.g. "0.7840" -> 0.784). The original style is
# restored immediately afterwards so no formatting changes leak in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $range = $ws.Range($CellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    $range.Style = $origStyle
}

# --- Price column (D) refreshes ---
Set-TextValue "D2"  "267.75"
Set-TextValue "D3"  "22.67"
Set-TextValue "D4"  "6.334"
Set-TextValue "D5"  "0.06164"
Set-TextValue "D6"  "3.645"
Set-TextValue "D7"  "6.664"
Set-TextValue "D8"  "1.383"
Set-TextValue "D9"  "0.8317"
Set-TextValue "D10" "0.01362"
Set-TextValue "D11" "0.1599"
Set-TextValue "D12" "0.08289"
Set-TextValue "D14" "0.03223"

# --- Rows 15 & 16 swap places (BitMartToken now ranks above MCDex) ---
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09306"
Set-TextValue "E15" "14BitMartTokenBMX"

Set-TextValue "B16" "MCDex"
Set-TextValue "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.814"
Set-TextValue "E16" "15MCDexMCB"

# --- Remaining price column (D) refreshes ---
Set-TextValue "D17" "0.001648"
Set-TextValue "D18" "0.04762"
Set-TextValue "D19" "0.006372"
Set-TextValue "D20" "0.005652"
Set-TextValue "D21" "0.001079"
Set-TextValue "D22" "0.0001505"
Set-TextValue "D23" "3.725"
Set-TextValue "D24" "2.413"
Set-TextValue "D25" "0.3333"
Set-TextValue "D27" "0.0002711"
Set-TextValue "D40" "0.04705"
Set-TextValue "D41" "0.006975"
Set-TextValue "D43" "0.003531"
Set-TextValue "D44" "0.01177"
Set-TextValue "D45" "0.00006274"
Set-TextValue "D46" "0.0009926"
Set-TextValue "D48" "0.7840"
Set-TextValue "D50" "0.00002406"
